{"js": "// Update the three-digit x one-digit multiplication problems/answers\n// inside the table cells. Each old expression is unique in the\n// document, so an exact (case-sensitive, whole-match) search safely\n// identifies the single run to replace.\nconst replacements = [\n  [\"809\u00d79=7281\", \"154\u00d72=308\"],\n  [\"778\u00d74=3112\", \"746\u00d75=3730\"],\n  [\"949\u00d74=3796\", \"782\u00d78=6256\"],\n  [\"111\u00d78=888\", \"576\u00d75=2880\"],\n  [\"109\u00d79=981\", \"912\u00d78=7296\"],\n  [\"877\u00d75=4385\", \"448\u00d75=2240\"],\n  [\"114\u00d77=798\", \"640\u00d77=4480\"],\n  [\"315\u00d79=2835\", \"922\u00d75=4610\"],\n  [\"982\u00d79=8838\", \"607\u00d74=2428\"],\n  [\"552\u00d76=3312\", \"967\u00d78=7736\"],\n  [\"786\u00d72=1572\", \"961\u00d75=4805\"],\n  [\"396\u00d78=3168\", \"843\u00d73=2529\"],\n  [\"442\u00d79=3978\", \"278\u00d73=834\"],\n  [\"456\u00d74=1824\", \"351\u00d74=1404\"],\n  [\"365\u00d74=1460\", \"627\u00d73=1881\"],\n  [\"855\u00d73=2565\", \"186\u00d76=1116\"],\n  [\"992\u00d72=1984\", \"748\u00d79=6732\"],\n  [\"619\u00d77=4333\", \"581\u00d79=5229\"],\n  [\"955\u00d72=1910\", \"401\u00d72=802\"],\n  [\"819\u00d76=4914\", \"746\u00d72=1492\"],\n  [\"912\u00d77=6384\", \"178\u00d76=1068\"],\n  [\"842\u00d73=2526\", \"243\u00d74=972\"],\n  [\"925\u00d74=3700\", \"597\u00d77=4179\"],\n  [\"313\u00d76=1878\", \"107\u00d75=535\"],\n  [\"908\u00d78=7264\", \"263\u00d77=1841\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the three-digit x one-digit multiplication problems/answers\n# inside the table cells. Each old expression is unique in the\n# document, so Find/Replace with MatchCase + MatchWholeWord safely\n# targets exactly one run per replacement.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"809\u00d79=7281\"; New = \"154\u00d72=308\" },\n    @{ Old = \"778\u00d74=3112\"; New = \"746\u00d75=3730\" },\n    @{ Old = \"949\u00d74=3796\"; New = \"782\u00d78=6256\" },\n    @{ Old = \"111\u00d78=888\";  New = \"576\u00d75=2880\" },\n    @{ Old = \"109\u00d79=981\";  New = \"912\u00d78=7296\" },\n    @{ Old = \"877\u00d75=4385\"; New = \"448\u00d75=2240\" },\n    @{ Old = \"114\u00d77=798\";  New = \"640\u00d77=4480\" },\n    @{ Old = \"315\u00d79=2835\"; New = \"922\u00d75=4610\" },\n    @{ Old = \"982\u00d79=8838\"; New = \"607\u00d74=2428\" },\n    @{ Old = \"552\u00d76=3312\"; New = \"967\u00d78=7736\" },\n    @{ Old = \"786\u00d72=1572\"; New = \"961\u00d75=4805\" },\n    @{ Old = \"396\u00d78=3168\"; New = \"843\u00d73=2529\" },\n    @{ Old = \"442\u00d79=3978\"; New = \"278\u00d73=834\" },\n    @{ Old = \"456\u00d74=1824\"; New = \"351\u00d74=1404\" },\n    @{ Old = \"365\u00d74=1460\"; New = \"627\u00d73=1881\" },\n    @{ Old = \"855\u00d73=2565\"; New = \"186\u00d76=1116\" },\n    @{ Old = \"992\u00d72=1984\"; New = \"748\u00d79=6732\" },\n    @{ Old = \"619\u00d77=4333\"; New = \"581\u00d79=5229\" },\n    @{ Old = \"955\u00d72=1910\"; New = \"401\u00d72=802\" },\n    @{ Old = \"819\u00d76=4914\"; New = \"746\u00d72=1492\" },\n    @{ Old = \"912\u00d77=6384\"; New = \"178\u00d76=1068\" },\n    @{ Old = \"842\u00d73=2526\"; New = \"243\u00d74=972\" },\n    @{ Old = \"925\u00d74=3700\"; New = \"597\u00d77=4179\" },\n    @{ Old = \"313\u00d76=1878\"; New = \"107\u00d75=535\" },\n    @{ Old = \"908\u00d78=7264\"; New = \"263\u00d77=1841\" }\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Old, $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, $pair.New, $wdReplaceAll)\n}\n"}
